# feat: compared higher dim AE on 2d PCA
#
# The "Eigenfaces: 3-4 components ... after PCA" note (cell B26) is extended
# with the content that used to live in D34 ("try AE with higher
# dimensionality and use PCA to visualize"), and the now-redundant D34 cell
# is cleared. Row 26 grows taller to fit the longer text, and the active
# selection ends up on B26 (the cell that was edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$eigenfacesNote = $ws.Range("B26").Value2
$aeNote = $ws.Range("D34").Value2

$ws.Range("B26").Value2 = $eigenfacesNote + ", " + $aeNote
$ws.Range("D34").Clear()

$ws.Rows.Item(26).RowHeight = 103

$ws.Range("B26").Select()
